# Apply crypto price/volume updates for Tue Aug  8 21:20:19 UTC 2023 run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to store plain text even when the string looks numeric
    # (e.g. "9.350", "0.5094"), matching the source inlineStr cells, then
    # restore the default "Normal" style so no stray formatting is left behind.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "30.007.39"
$ws.Range("E2").Value = "  +2.95%  "

Set-TextValue $ws.Range("D3") "1.865.82"
$ws.Range("E3").Value = "  +2.25%  "

Set-TextValue $ws.Range("D4") "0.9993"
$ws.Range("E4").Value = "  -0.06%  "

Set-TextValue $ws.Range("D5") "246.16"
$ws.Range("E5").Value = "  +1.92%  "

Set-TextValue $ws.Range("D6") "0.6406"
$ws.Range("E6").Value = "  +3.84%  "

Set-TextValue $ws.Range("D7") "0.9997"
$ws.Range("E7").Value = "  -0.03%  "

Set-TextValue $ws.Range("D8") "0.3004"
$ws.Range("E8").Value = "  +3.99%  "

Set-TextValue $ws.Range("D9") "0.07501"
$ws.Range("E9").Value = "  +2.25%  "

Set-TextValue $ws.Range("D10") "24.41"
$ws.Range("E10").Value = "  +6.40%  "

Set-TextValue $ws.Range("D11") "0.07687"
$ws.Range("E11").Value = "  +0.14%  "

Set-TextValue $ws.Range("D12") "1.874.28"
$ws.Range("E12").Value = "  +2.84%  "

Set-TextValue $ws.Range("D13") "5.072"
$ws.Range("E13").Value = "  +2.35%  "

$ws.Range("E14").Value = "  +4.56%  "

Set-TextValue $ws.Range("D15") "84.39"
$ws.Range("E15").Value = "  +2.95%  "

Set-TextValue $ws.Range("D16") "0.000009473"
$ws.Range("E16").Value = "  +6.45%  "

Set-TextValue $ws.Range("D17") "6.110"
$ws.Range("E17").Value = "  +4.76%  "

Set-TextValue $ws.Range("D18") "29.978.49"
$ws.Range("E18").Value = "  +2.93%  "

Set-TextValue $ws.Range("D19") "2.122.24"
$ws.Range("E19").Value = "  +2.67%  "

Set-TextValue $ws.Range("D20") "241.90"
$ws.Range("E20").Value = "  +1.74%  "

$ws.Range("E21").Value = "  +2.42%  "

Set-TextValue $ws.Range("D22") "0.9996"
$ws.Range("E22").Value = "  -0.09%  "

Set-TextValue $ws.Range("D23") "7.479"
$ws.Range("E23").Value = "  +4.46%  "

Set-TextValue $ws.Range("D24") "1.001"
$ws.Range("E24").Value = "  -0.07%  "

Set-TextValue $ws.Range("D25") "159.67"
$ws.Range("E25").Value = "  +1.08%  "

Set-TextValue $ws.Range("D26") "0.1426"
$ws.Range("E26").Value = "  +0.72%  "

Set-TextValue $ws.Range("D27") "8.590"
$ws.Range("E27").Value = "  +1.99%  "

Set-TextValue $ws.Range("D28") "18.07"
$ws.Range("E28").Value = "  +2.52%  "

Set-TextValue $ws.Range("D29") "0.06134"
$ws.Range("E29").Value = "  +10.46%  "

Set-TextValue $ws.Range("D30") "1.505"
$ws.Range("E30").Value = "  +1.38%  "

$ws.Range("E31").Value = "  +5.81%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D32") "4.155"
$ws.Range("E32").Value = "  +1.47%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D33") "4.135"
$ws.Range("E33").Value = "  +1.08%  "

Set-TextValue $ws.Range("D34") "1.871"
$ws.Range("E34").Value = "  +2.80%  "

Set-TextValue $ws.Range("D35") "1.168"

Set-TextValue $ws.Range("D36") "0.7311"
$ws.Range("E36").Value = "  -0.62%  "

Set-TextValue $ws.Range("D37") "2.604"
$ws.Range("E37").Value = "  -0.35%  "

Set-TextValue $ws.Range("D38") "2.869"
$ws.Range("E38").Value = "  +1.12%  "

Set-TextValue $ws.Range("D39") "0.01804"
$ws.Range("E39").Value = "  +2.75%  "

Set-TextValue $ws.Range("D40") "1.223.61"

$ws.Range("E41").Value = "  +1.13%  "

Set-TextValue $ws.Range("D42") "6.280"
$ws.Range("E42").Value = "  -0.56%  "

Set-TextValue $ws.Range("D43") "1.002"
$ws.Range("E43").Value = "  +0.14%  "

Set-TextValue $ws.Range("D44") "2.029.32"
$ws.Range("E44").Value = "  +2.95%  "

Set-TextValue $ws.Range("D45") "102.29"
$ws.Range("E45").Value = "  +0.93%  "

Set-TextValue $ws.Range("D46") "66.59"
$ws.Range("E46").Value = "  +3.03%  "

$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D47") "0.00000000122"
$ws.Range("E47").Value = "  -5.79%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D48") "0.5094"
$ws.Range("E48").Value = "  +0.16%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D49") "9.350"
$ws.Range("E49").Value = "  +3.34%  "

$ws.Range("E50").Value = "  +2.52%  "

$ws.Range("E51").Value = "  +3.05%  "
